{"js": "// Change the \"CoursePlus (with Zoom information!):\" text to the plain\n// \"CoursePlus:\" label, and add a new bold paragraph right after the\n// \"Class Website / CoursePlus\" paragraph announcing that the Zoom link\n// will be emailed to students.\n\nconst body = context.document.body;\n\n// Locate the paragraph that currently holds the CoursePlus/Zoom text so we\n// can anchor the new paragraph to it (do this before any text edit so the\n// paragraph's own formatting/style is still pristine).\nconst searchResults = body.search(\"CoursePlus (with Zoom information!):\", {\n  matchCase: true,\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"CoursePlus (with Zoom information!):\" in the document.');\n}\n\nconst coursePlusRange = searchResults.items[0];\nconst coursePlusParagraph = coursePlusRange.paragraphs.getFirst();\n\n// Insert the new \"Zoom link will be emailed to students.\" paragraph right\n// after the CoursePlus paragraph, inheriting its (BodyText) style, then make\n// its text bold.\nconst zoomParagraph = coursePlusParagraph.insertParagraph(\n  \"Zoom link will be emailed to students.\",\n  Word.InsertLocation.after\n);\nzoomParagraph.font.bold = true;\n\n// Now collapse the \"CoursePlus (\", \"with Zoom information!\", \"):\" runs down\n// to a single plain \"CoursePlus:\" run.\ncoursePlusRange.insertText(\"CoursePlus:\", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Change the \"CoursePlus (with Zoom information!):\" text to the plain\n# \"CoursePlus:\" label, and add a new bold paragraph right after the\n# \"Class Website / CoursePlus\" paragraph announcing that the Zoom link\n# will be emailed to students.\n\n$doc = $word.ActiveDocument\n\n# Locate the paragraph that currently holds \"CoursePlus (with Zoom\n# information!):\" (it's the same paragraph as \"Class Website:\").\n$target = $doc.Content\n$target.Find.Execute(\"CoursePlus (with Zoom information!):\") | Out-Null\n$coursePlusParaIndex = $target.Paragraphs(1).Index\n\n$coursePlusPara = $doc.Paragraphs($coursePlusParaIndex)\n\n# Step 1: insert the new \"Zoom link will be emailed to students.\" paragraph\n# right after it, appending at the clean paragraph end so none of the\n# existing runs are disturbed.\n$endPos = $coursePlusPara.Range.End\n$insertPoint = $doc.Range($endPos, $endPos)\n$insertPoint.InsertAfter(\"Zoom link will be emailed to students.`r\")\n\n$zoomPara = $doc.Paragraphs($coursePlusParaIndex + 1)\n$zoomStart = $zoomPara.Range.Start\n$zoomEnd = $zoomPara.Range.End\n$zoomTextRange = $doc.Range($zoomStart, $zoomEnd - 1)\n$zoomTextRange.Font.Bold = 1\n\n# Step 2: collapse the \"CoursePlus (\", \"with Zoom information!\", \"):\" runs\n# down to a single plain \"CoursePlus:\" run.\n$findRange = $doc.Content\n$findRange.Find.Execute(\"CoursePlus (with Zoom information!):\") | Out-Null\n$findRange.Text = \"CoursePlus:\"\n"}
